$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.13"
$ws.Range("E2").Value = "'-1.28%"
$ws.Range("D3").Value = "'48.62"
$ws.Range("E3").Value = "'8.66%"
$ws.Range("D4").Value = "'5.243"
$ws.Range("E4").Value = "'2.01%"
$ws.Range("D5").Value = "'0.07788"
$ws.Range("E5").Value = "'-3.41%"
$ws.Range("D6").Value = "'4.514"
$ws.Range("E6").Value = "'-0.15%"
$ws.Range("D7").Value = "'1.296"
$ws.Range("E7").Value = "'19.06%"
$ws.Range("D8").Value = "'1.562"
$ws.Range("E8").Value = "'-7.71%"
$ws.Range("E9").Value = "'-3.85%"
$ws.Range("D10").Value = "'0.1928"
$ws.Range("E10").Value = "'0.44%"
$ws.Range("D11").Value = "'0.09272"
$ws.Range("E11").Value = "'-1.33%"
$ws.Range("E12").Value = "'7.11%"
$ws.Range("D13").Value = "'0.1048"
$ws.Range("E13").Value = "'0.41%"
$ws.Range("D14").Value = "'0.001299"
$ws.Range("E14").Value = "'-1.44%"
$ws.Range("D15").Value = "'0.04201"
$ws.Range("E15").Value = "'-0.69%"
$ws.Range("D16").Value = "'0.005832"
$ws.Range("E16").Value = "'-1.95%"
$ws.Range("E17").Value = "'-1.92%"
$ws.Range("D18").Value = "'2.406"
$ws.Range("E18").Value = "'-0.21%"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'1.49%"
$ws.Range("D20").Value = "'8.138"
$ws.Range("E20").Value = "'-1.90%"
$ws.Range("D21").Value = "'0.1367"
$ws.Range("E21").Value = "'-1.22%"
$ws.Range("D22").Value = "'0.3080"
$ws.Range("E22").Value = "'-1.99%"
$ws.Range("E23").Value = "'1.77%"
$ws.Range("D24").Value = "'0.004222"
$ws.Range("E24").Value = "'-6.69%"
$ws.Range("D25").Value = "'0.0001361"
$ws.Range("E25").Value = "'1.41%"
$ws.Range("D38").Value = "'0.02562"
$ws.Range("E38").Value = "'-4.62%"
$ws.Range("D39").Value = "'0.05800"
$ws.Range("E39").Value = "'6.32%"
$ws.Range("D40").Value = "'0.01082"
$ws.Range("E40").Value = "'97.25%"
$ws.Range("D41").Value = "'0.008138"
$ws.Range("E41").Value = "'5.22%"
$ws.Range("D42").Value = "'0.1428"
$ws.Range("E42").Value = "'0.57%"
$ws.Range("D43").Value = "'0.008457"
$ws.Range("E43").Value = "'15.29%"
$ws.Range("D44").Value = "'0.008524"
$ws.Range("E44").Value = "'-0.66%"
$ws.Range("D45").Value = "'0.3118"
$ws.Range("E45").Value = "'-0.58%"
$ws.Range("D46").Value = "'0.00006914"
$ws.Range("E46").Value = "'1.64%"
$ws.Range("E47").Value = "'1.25%"
$ws.Range("D48").Value = "'0.05508"
$ws.Range("E48").Value = "'-20.60%"
$ws.Range("D49").Value = "'0.004032"
$ws.Range("E49").Value = "'1.25%"
$ws.Range("E50").Value = "'1.25%"
$ws.Range("E51").Value = "'1.25%"
